$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two review-confirmation flags were changed from "confirm" to "no"
$ws.Range("G19").Value = "no"
$ws.Range("G24").Value = "no"

# Update the saved selection to match the author's final cursor position
$ws.Range("G25").Select()
